$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Voeg de tijdsbesteding toe voor de week van maandag 19/10/2015
$ws.Range("B4").Value = 3

# Selectie verplaatsen naar B5 (zoals na het invoeren van een waarde in B4)
$ws.Range("B5").Select()
